# Add the missing "Bedtime" variable definition to the Dictionary sheet.
#
# A new data row is inserted right above the existing "SleepInterval" row
# (currently row 41), pushing it and everything below it down by one row.
# The new row documents the "Bedtime" variable (name / short description /
# long description). The existing "SleepInBed" row's long description
# (which explains it is the same as "Sleep") is also reworded to clarify
# that note has applied "since ActiPASS version 1.50".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Dictionary")

# Insert a new blank row at position 41 (shifts SleepInterval etc. down by one).
$ws.Rows.Item(41).Insert()

# Populate the new "Bedtime" variable-definition row.
$ws.Range("A41").Value = "Bedtime"
$ws.Range("B41").Value = "The time within this calendar day which is flagged as belonging to a bedtime"

# Update the wording of the "SleepInBed" row's long description (now at
# row 44, column C, after the insert above) to reference ActiPASS 1.50.
$ws.Range("C44").Value = 'This is same as "Sleep" now since ActiPASS version 1.50. Possible sleep outside bedtimes are now named "LieStill"'

# Long description for the new "Bedtime" row.
$ws.Range("C41").Value = "A calendar day usually contains more than one (full or partial) bedtimes. This variable is calculated by counting the number of epochs (1s) which are flagged as belonging to a bedtime. "
